$p = $ppt.ActivePresentation

$s2 = $p.Slides.Item(2)
$notes2 = $s2.NotesPage
$notes2.Shapes.Item(2).TextFrame.TextRange.Text = "Some notes on the second slide."

$s3 = $p.Slides.Item(3)
$notes3 = $s3.NotesPage
$notes3.Shapes.Item(2).TextFrame.TextRange.Text = "Final notes on the third slide.`rSecond line of notes."
